# Update "Forecast Comparison" sheet: shift weekly data forward by one week
# (new row 17 introduces the 2025-05-18 week), per corrected forecast export.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

$ws.Range("B2").Value = "'2025-02-02"
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 2
$ws.Range("G2").Value = 5
$ws.Range("H2").Value = 11

$ws.Range("B3").Value = "'2025-02-09"
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 2
$ws.Range("G3").Value = 5
$ws.Range("H3").Value = 10

$ws.Range("B4").Value = "'2025-02-16"
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 5
$ws.Range("F4").Value = 4
$ws.Range("G4").Value = 8
$ws.Range("H4").Value = 16

$ws.Range("B5").Value = "'2025-02-23"
$ws.Range("D5").Value = 2
$ws.Range("E5").Value = 5
$ws.Range("F5").Value = 4
$ws.Range("G5").Value = 8
$ws.Range("H5").Value = 16

$ws.Range("B6").Value = "'2025-03-02"
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 6
$ws.Range("F6").Value = 6
$ws.Range("G6").Value = 10
$ws.Range("H6").Value = 18

$ws.Range("B7").Value = "'2025-03-09"
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = 6
$ws.Range("F7").Value = 6
$ws.Range("G7").Value = 10
$ws.Range("H7").Value = 19

$ws.Range("B8").Value = "'2025-03-16"
$ws.Range("D8").Value = 1
$ws.Range("E8").Value = 6
$ws.Range("F8").Value = 6
$ws.Range("G8").Value = 11
$ws.Range("H8").Value = 20

$ws.Range("B9").Value = "'2025-03-23"
$ws.Range("D9").Value = 1
$ws.Range("E9").Value = 6
$ws.Range("F9").Value = 6
$ws.Range("G9").Value = 11
$ws.Range("H9").Value = 19

$ws.Range("B10").Value = "'2025-03-30"
$ws.Range("D10").Value = 1
$ws.Range("E10").Value = 6
$ws.Range("F10").Value = 5
$ws.Range("G10").Value = 9
$ws.Range("H10").Value = 18

$ws.Range("B11").Value = "'2025-04-06"
$ws.Range("D11").Value = 1
$ws.Range("E11").Value = 6
$ws.Range("F11").Value = 5
$ws.Range("G11").Value = 10
$ws.Range("H11").Value = 19

$ws.Range("B12").Value = "'2025-04-13"
$ws.Range("D12").Value = 2
$ws.Range("E12").Value = 7
$ws.Range("F12").Value = 6
$ws.Range("G12").Value = 11
$ws.Range("H12").Value = 20

$ws.Range("B13").Value = "'2025-04-20"
$ws.Range("D13").Value = 2
$ws.Range("E13").Value = 7
$ws.Range("F13").Value = 6
$ws.Range("G13").Value = 11
$ws.Range("H13").Value = 20

$ws.Range("B14").Value = "'2025-04-27"
$ws.Range("D14").Value = 2
$ws.Range("E14").Value = 6
$ws.Range("F14").Value = 6
$ws.Range("G14").Value = 10
$ws.Range("H14").Value = 20

$ws.Range("B15").Value = "'2025-05-04"
$ws.Range("D15").Value = 2
$ws.Range("E15").Value = 6
$ws.Range("F15").Value = 4
$ws.Range("G15").Value = 9
$ws.Range("H15").Value = 19

$ws.Range("B16").Value = "'2025-05-11"
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 6
$ws.Range("F16").Value = 4
$ws.Range("G16").Value = 9
$ws.Range("H16").Value = 18

$ws.Range("B17").Value = "'2025-05-18"
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 6
$ws.Range("F17").Value = 4
$ws.Range("G17").Value = 8
$ws.Range("H17").Value = 18

# Update "Summary" sheet to reflect the corrected historical range and forecast stats
$ws2 = $wb.Worksheets.Item("Summary")
$ws2.Range("B2").Value = "2022-12-25 to 2025-01-26"
$ws2.Range("B8").Value = "53 units"
$ws2.Range("B9").Value = "'22"
$ws2.Range("B10").Value = "'9"
$ws2.Range("B11").Value = "'4"
$ws2.Range("B13").Value = "'2025-04-20"
$ws2.Range("B15").Value = "'2025-02-02"
